$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 138, which shifts existing rows 138-263 down to 139-264
$ws.Rows.Item(138).Insert()

# Fill in the new row 138 with the new data record
$ws.Cells.Item(138, 1).Value = 8
$ws.Cells.Item(138, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(138, 3).Value = "Coquimbo"
$ws.Cells.Item(138, 4).Value = 44904
$ws.Cells.Item(138, 5).Value = 4
$ws.Cells.Item(138, 6).Value = 100112037
$ws.Cells.Item(138, 7).Value = "Cebollín"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 1000
$ws.Cells.Item(138, 11).Value = 1200
$ws.Cells.Item(138, 12).Value = 1400
$ws.Cells.Item(138, 13).Value = 1300
$ws.Cells.Item(138, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(138, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(138, 16).Value = 217
$ws.Cells.Item(138, 17).Value = 6
$ws.Cells.Item(138, 18).Value = "Hortaliza"
